$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("period_lbl")
$ws2 = $wb.Worksheets.Item("time_period_type")
$ws3 = $wb.Worksheets.Item("year")

# --- Rename the "label" columns first (B1 on each sheet) ---
$ws1.Range("B1").Value = "period_lbl"
$ws2.Range("B1").Value = "time_period_type"
$ws3.Range("B1").Value = "year"

# --- Then rename the "code" columns (A1 on each sheet) ---
$ws1.Range("A1").Value = "period_code"
$ws2.Range("A1").Value = "time_period_code"
$ws3.Range("A1").Value = "year_code"

# --- Apply the small monospace "code-ish" font style to the renamed label columns ---
$f1 = $ws1.Range("B1").Font
$f1.Name = "Consolas"
$f1.Size = 7
$f1.Color = 7901646
$ws1.Range("B1").VerticalAlignment = -4108

$f2 = $ws2.Range("B1").Font
$f2.Name = "Consolas"
$f2.Size = 7
$f2.Color = 7901646
$ws2.Range("B1").VerticalAlignment = -4108

# --- Page setup on "period_lbl" ---
$ps1 = $ws1.PageSetup
$ps1.PaperSize = 9
$ps1.Orientation = 1

# --- View / selection state: land on sheet "year" as the active tab/sheet ---
[void]$ws1.Range("A1").Select()

$ws2.Activate()
[void]$ws2.Range("D8").Select()

$ws3.Activate()
[void]$ws3.Range("D14").Select()
